$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the NZSCED field-of-study reference table:
# Code / Description / Four_Digit_Code / Four_Digit_Description /
# Two_Digit_Code / Two_Digit_Description / Definition
$rows = @(
    @{ Code="444444"; Desc="Don't Know";               Code4="4444"; Code2="44" },
    @{ Code="555555"; Desc="Refused to Answer";         Code4="5555"; Code2="55" },
    @{ Code="777777"; Desc="Response Unidentifiable";   Code4="7777"; Code2="77" },
    @{ Code="888888"; Desc="Response Outside Scope";    Code4="8888"; Code2="88" },
    @{ Code="999999"; Desc="Not Stated";                Code4="9999"; Code2="99" }
)

$startRow = 464
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $definition = "$($row.Desc)/$($row.Desc)/$($row.Desc)"

    # Columns A, C, E hold codes that look like plain numbers (e.g. "444444").
    # Force them to be stored as text (matching the source data) rather than
    # letting them be auto-converted to numeric values, then restore the
    # cell style back to Normal so no extra formatting is left behind.
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row.Code
    $ws.Range("A$r").Style = "Normal"

    $ws.Range("B$r").Value = $row.Desc

    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $row.Code4
    $ws.Range("C$r").Style = "Normal"

    $ws.Range("D$r").Value = $row.Desc

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row.Code2
    $ws.Range("E$r").Style = "Normal"

    $ws.Range("F$r").Value = $row.Desc

    $ws.Range("G$r").Value = $definition
}
